$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rows 73 and 74 had their match data swapped (id/div/div-name/date
#    stay put, everything else - match id, teams, score, odds - moves
#    from one row to the other).
# ------------------------------------------------------------------
$swapCols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")
foreach ($col in $swapCols) {
    $cellA = $col + "73"
    $cellB = $col + "74"
    $valA = $ws.Range($cellA).Value2
    $valB = $ws.Range($cellB).Value2
    $ws.Range($cellA).Value2 = $valB
    $ws.Range($cellB).Value2 = $valA
}

# ------------------------------------------------------------------
# 2) Six new upcoming fixtures were appended as rows 132-137.
# ------------------------------------------------------------------
$newRows = @(
    @{ Row=132; A=130; B=7127395; E=45380.23958333334; F="Adelaide United";     G="Western United FC";       K=1.666; L=4.2;   M=4.333; N=1.85;  O=4.5;  P=3.4;  Q=-0.5;  R=1.89; S=2.01; T=3.5;  U=2;     V=1.85;  W=0; X=0; Y=0; Z=0; AA=0 },
    @{ Row=133; A=131; B=7126793; E=45381.14583333334; F="Melbourne City";       G="Newcastle Jets";          K=1.571; L=4.333; M=5;     N=1.533; O=4.75; P=5.25; Q=-1;    R=1.9;  S=2;    T=3.25; U=1.875; V=1.975; W=0; X=0; Y=0; Z=0; AA=0 },
    @{ Row=134; A=132; B=7127396; E=45381.23958333334; F="Sydney FC";            G="Central Coast Mariners";  K=2.15;  L=3.6;   M=3.1;   N=2.15;  O=3.75; P=3.1;  Q=-0.25; R=1.93; S=1.97; T=2.75; U=1.8;   V=2.05;  W=0; X=0; Y=0; Z=0; AA=0 },
    @{ Row=135; A=133; B=7127394; E=45381.875;          F="Wellington Phoenix";  G="Brisbane Roar";           K=1.8;   L=3.8;   M=4;     N=2;     O=3.6;  P=3.6;  Q=-0.5;  R=2.06; S=1.84; T=3;    U=1.925; V=1.925; W=0; X=0; Y=0; Z=0; AA=0 },
    @{ Row=136; A=134; B=7127397; E=45382.04166666666; F="Melbourne Victory";    G="Perth Glory";             K=1.4;   L=5;     M=6.5;   N=1.4;   O=5.5;  P=6.5;  Q=-1.25; R=1.91; S=1.99; T=3.25; U=1.925; V=1.925; W=0; X=0; Y=0; Z=0; AA=0 },
    @{ Row=137; A=135; B=7127398; E=45383.04166666666; F="Macarthur FC";         G="Western Sydney Wanderers";K=2.5;   L=3.5;   M=2.625; N=2.625; O=3.75; P=2.45; Q=0;     R=2.03; S=1.87; T=3.25; U=1.95;  V=1.9;   W=0; X=0; Y=0; Z=0; AA=0 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = "Australia ALeague"
    $ws.Cells.Item($row, 4).Value = "Australia ALeague"
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G

    # Match the look of the rest of the table: column A is bold/centered
    # with a thin border (same style as the other "id" cells), column E
    # uses the custom date/time number format.
    $ws.Range("A131").Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws.Range("E131").Copy()
    $ws.Cells.Item($row, 5).PasteSpecial(-4122)

    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
    $ws.Cells.Item($row, 21).Value = $r.U
    $ws.Cells.Item($row, 22).Value = $r.V
    $ws.Cells.Item($row, 23).Value = $r.W
    $ws.Cells.Item($row, 24).Value = $r.X
    $ws.Cells.Item($row, 25).Value = $r.Y
    $ws.Cells.Item($row, 26).Value = $r.Z
    $ws.Cells.Item($row, 27).Value = $r.AA
}
